# Bump the "Version: 1.0.2" text to "Version: 1.0.3" and add a new
# translation row ("SingleUseId56" / "PCB: <value> °<value>") on the
# Translation sheet, per commit "Updated source code to version 1.0.3".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translation")

# Row 25, column F holds the firmware version string shown on screen.
$ws.Cells.Item(25, 6).Value = "Version: 1.0.3"

# Row 59 was a blank placeholder row; fill it in with the new PCB
# temperature text entry (TEXT ID / TYPOGRAPHY NAME / ALIGNMENT /
# DIRECTION / GB columns, mirroring the existing rows' layout).
$ws.Cells.Item(59, 2).Value = "SingleUseId56"
$ws.Cells.Item(59, 3).Value = "text_graph_size"
$ws.Cells.Item(59, 4).Value = "Left"
$ws.Cells.Item(59, 5).Value = "LTR"
$ws.Cells.Item(59, 6).Value = "PCB: <value> °<value>        "
